$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.541.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.390.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'407.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'126.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.610"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.55%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.708"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.84%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "  -8.87%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'41.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.926.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'8.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0000203"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.97%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'20.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.53%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.461.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'12.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.00%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'61.712.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'479.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +20.73%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'88.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.60%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'13.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'3.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'32.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'8.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.65%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'11.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.164"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.47%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.111"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.48%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'40.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.84%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'56.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.03%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.0476"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.36%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'149.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.10%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "  -2.32%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "TheGraph"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.313"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "Stacks"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'2.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "  +2.24%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'2.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.06%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'4.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'2.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +16.24%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "Celestia"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'16.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.21%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'21.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.141"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.06%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'111.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.30%  "
$ws.Range("E51").Style = "Normal"
